$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column A. The engine round-trips "characters" width through a pixel
# grid using the default font metrics, so 19.15 in is the value that lands on
# a stored width of exactly 20 (the target value) after that conversion.
$ws.Columns("A").ColumnWidth = 19.15

# Copy the formatting (font/style) of the last existing data row down onto the
# new rows, then fill in the new values.
$ws.Range("A5:B5").Copy()
[void]$ws.Range("A6:B14").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Column A gets filled top-to-bottom first, then column B - this matches the
# order new strings were appended to the shared-string table in the target.
$colA = @("Teh chamomile", "Teh jahe", "Teh melati ", "Teh peppermint ", "Teh herbal ", "Teh kembang sepatu", "Teh merah", "Teh hibiscus", "Teh Pu-erh")
$colB = @("Artha Cipta", "Bahana", "Citalahab", "Baru Ulis", "Berdikari", "Bodjong", "Bintang", "Mulia", "Bolindo")

$r = 6
foreach ($v in $colA) {
    $ws.Cells.Item($r, 1).Value = $v
    $r = $r + 1
}

$r = 6
foreach ($v in $colB) {
    $ws.Cells.Item($r, 2).Value = $v
    $r = $r + 1
}

$ws.Rows("6:14").RowHeight = 15.5

# Move the active selection to C5, matching the saved workbook view.
[void]$ws.Range("C5").Select()
